$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new product row (Leche ultrapasteurizada, Manfrey) was inserted as row 4 in the
# source workbook, pushing every existing data row (old rows 4-35) down by one
# (new rows 5-36). Reproduce that by copying each existing row one position down,
# working from the bottom up so we never overwrite a row before it has been copied.
for ($r = 35; $r -ge 4; $r--) {
    $src = $ws.Range("A" + $r + ":O" + $r)
    $dst = $ws.Range("A" + ($r + 1) + ":O" + ($r + 1))
    $src.Copy($dst)
}

# Row 4 currently holds a duplicate of the old row 4 (now also in row 5); wipe it so
# we can write the brand-new record cleanly (this keeps the existing A4/O4 styles).
$ws.Range("A4:O4").ClearContents()

# New row 4: Leche / ultrapasteurizada / entera / Manfrey, barcode 7791058000595.
# Only Codigo, TipoArticulo, Descripcion, Variedad, Marca, Pesable,
# TieneVencimiento and ImagenExactaDelArticulo are populated.
$ws.Range("A4").Value = 7791058000595
$ws.Range("B4").Value = "Leche"
$ws.Range("C4").Value = "ultrapasteurizada"
$ws.Range("D4").Value = "entera"
$ws.Range("E4").Value = "Manfrey"
$ws.Range("L4").Value = $false
$ws.Range("M4").Value = $false
$ws.Range("O4").Value = $false

# Row 5 (the old row 4 - Salame tipo milan Fela) gets its Descripcion updated.
$ws.Range("C5").Value = "pelado"
